$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.55"
$ws.Range("G2").Value = "'3"
$ws.Range("D3").Value = "'23.26"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'5.409"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.05976"
$ws.Range("G5").Value = "'3"
$ws.Range("D6").Value = "'3.437"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'6.538"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'0.8104"
$ws.Range("G8").Value = "'3"
$ws.Range("D9").Value = "'0.9340"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.1425"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.07438"
$ws.Range("G11").Value = "'3"
$ws.Range("D12").Value = "'0.03301"
$ws.Range("G12").Value = "'3"
$ws.Range("D13").Value = "'0.03073"
$ws.Range("G13").Value = "'3"
$ws.Range("D14").Value = "'0.09365"
$ws.Range("G14").Value = "'3"
$ws.Range("D15").Value = "'3.857"
$ws.Range("G15").Value = "'3"
$ws.Range("D16").Value = "'0.001578"
$ws.Range("G16").Value = "'3"
$ws.Range("D17").Value = "'0.04706"
$ws.Range("G17").Value = "'3"
$ws.Range("D18").Value = "'0.0005923"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "'3"
$ws.Range("D19").Value = "'0.005928"
$ws.Range("G19").Value = "'3"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.004904"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'0.00006804"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'3.572"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'2.151"
$ws.Range("G24").Value = "'3"
$ws.Range("G25").Value = "'3"
$ws.Range("G26").Value = "'3"
$ws.Range("D27").Value = "'0.0002341"
$ws.Range("G27").Value = "'3"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("G38").Value = "'3"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.03975"
$ws.Range("G40").Value = "'3"
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D41").Value = "'0.005103"
$ws.Range("E41").Value = "40CEJICEJIBestin24h"
$ws.Range("G41").Value = "'3"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006394"
$ws.Range("E42").Value = "41KickTokenKICK"
$ws.Range("G42").Value = "'3"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1077"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.009212"
$ws.Range("G44").Value = "'3"
$ws.Range("D45").Value = "'0.00005229"
$ws.Range("G45").Value = "'3"
$ws.Range("G46").Value = "'3"
$ws.Range("D47").Value = "'0.7254"
$ws.Range("G47").Value = "'3"
$ws.Range("D48").Value = "'0.002409"
$ws.Range("G48").Value = "'3"
$ws.Range("G49").Value = "'3"
$ws.Range("G50").Value = "'3"
$ws.Range("G51").Value = "'3"
